$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Rigorous(tough):: Government should make rigorous law to control crime."
$ws.Range("A11").Value = "Consent(permission):: Without taking any consent from head of the company, the manager took th ewrong action"
$ws.Range("A12").Value = "Hypothesis(assumption)::Engineer wrong Hypothesis causes weak construction"
$ws.Range("A13").Value = "Peculiar(unique)::Taj Mahal is the peculiar architecture in the world."
$ws.Range("A14").Value = "Adequate(enough)::In nearly every paper there were adequate extracts."

$ws.Columns.Item(1).ColumnWidth = 14.45

$ws.Range("A14").Select()
